$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.657.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.963.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("E10").Value = "  -6.46%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.835"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.250.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("E16").Value = "  +3.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.968.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.630.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0856"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +5.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("E26").Value = "  +9.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +17.97%  "
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0616"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.99%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0212"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.372.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.141.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
